# "Alguns conflitos na aula 6"
# Fix the column header in L3 (foreign-key column under the "Email" block)
# from the wrong value "Pessoa" to the correct value "idPessoa", matching
# the pattern used by the other relationship tables (e.g. H3).
# Also update the worksheet's current selection to K13 (where the author
# was last working) to match the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = "idPessoa"

$ws.Range("K13").Select()
